# Auto commit at 2026-01-13  7:52:33.90
# Append the next day's (2026-01-12, serial 46034) two station rows to the
# bottom of the daily charging-data table on Sheet1 (rows 24 and 25),
# reusing the formatting of the last existing row (row 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 23's formatting down into the two new rows first, so the new
# cells pick up the same number formats/styles (date style, 0.00 style,
# integer style) as the rest of the table.
$ws.Range("A23:F23").Copy()
$ws.Range("A24:F24").PasteSpecial()
$ws.Range("A23:F23").Copy()
$ws.Range("A25:F25").PasteSpecial()

# Row 24: 四方坪站 (station 1)
$ws.Cells.Item(24, 1).Value = 46034
$ws.Cells.Item(24, 2).Value = "四方坪站"
$ws.Cells.Item(24, 3).Value = 14159.88
$ws.Cells.Item(24, 4).Value = 9722.2099999999991
$ws.Cells.Item(24, 5).Value = 3186.48
$ws.Cells.Item(24, 6).Value = 609

# Row 25: 高岭站 (station 2)
$ws.Cells.Item(25, 1).Value = 46034
$ws.Cells.Item(25, 2).Value = "高岭站"
$ws.Cells.Item(25, 3).Value = 5338.08
$ws.Cells.Item(25, 4).Value = 4389.1000000000004
$ws.Cells.Item(25, 5).Value = 1409.82
$ws.Cells.Item(25, 6).Value = 176

# Match the workbook's saved view state (selection / active cell) from the diff.
[void]$ws.Range("H24").Select()
